$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A9: 白/球比例 -> 白/球比值
$ws.Range("A9").Value = "白/球比值"

# E10, E11, E12: umol/L -> μmol/L
$ws.Range("E10").Value = "μmol/L"
$ws.Range("E11").Value = "μmol/L"
$ws.Range("E12").Value = "μmol/L"

# B13: HB sAg -> HBsAg
$ws.Range("B13").Value = "HBsAg"
